$d = $word.ActiveDocument

# Locate the paragraph that still contains the legacy Word field
# (begin fldChar / instrText runs / end fldChar) built from the
# "m:'doc.html'.fromHTMLURI()" instruction text, and rewrite it as a
# sequence of literal text runs (TokenIteratorFieldRewriterSplit style)
# instead of a live field, while keeping the _GoBack bookmark in place.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        foreach ($f in $p.Range.Fields) {
            if ($f.Code.Text -match "fromHTMLURI") {
                $target = $p
            }
        }
    }
}

if ($target -eq $null) {
    throw "Could not locate the fromHTMLURI() field paragraph"
}

$newXml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t>{</w:t></w:r>
<w:r><w:t>m</w:t></w:r>
<w:r><w:t>:</w:t></w:r>
<w:r><w:t>'</w:t></w:r>
<w:r><w:t>doc.html</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:t>'.fromHTMLURI()</w:t></w:r>
<w:r><w:t xml:space="preserve">}</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$target.Range.InsertXML($newXml)
